$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 257 (shifts existing rows 257:369 down to 260:372)
$ws.Range("A257:T259").EntireRow.Insert()

# --- New row 257 ---
$ws.Range("A257").Value = 3
$ws.Range("B257").Value = "Femacal de La Calera"
$ws.Range("C257").Value = "Coquimbo"
$ws.Range("D257").Value = 44875
$ws.Range("E257").Value = 5
$ws.Range("F257").Value = "Fruta"
$ws.Range("G257").Value = 100101
$ws.Range("H257").Value = "Berries"
$ws.Range("I257").Value = 100112025
$ws.Range("J257").Value = "Frutilla"
$ws.Range("K257").Value = "Sin especificar"
$ws.Range("L257").Value = "Especial"
$ws.Range("M257").Value = 75
$ws.Range("N257").Value = 8000
$ws.Range("O257").Value = 8000
$ws.Range("P257").Value = 8000
$ws.Range("Q257").Value = "$/bandeja 7 kilos"
$ws.Range("R257").Value = "Provincia de Melipilla"
$ws.Range("S257").Value = 1143
$ws.Range("T257").Value = 7

# --- New row 258 ---
$ws.Range("A258").Value = 3
$ws.Range("B258").Value = "Femacal de La Calera"
$ws.Range("C258").Value = "Coquimbo"
$ws.Range("D258").Value = 44875
$ws.Range("E258").Value = 5
$ws.Range("F258").Value = "Fruta"
$ws.Range("G258").Value = 100101
$ws.Range("H258").Value = "Berries"
$ws.Range("I258").Value = 100112025
$ws.Range("J258").Value = "Frutilla"
$ws.Range("K258").Value = "Sin especificar"
$ws.Range("L258").Value = "Primera"
$ws.Range("M258").Value = 78
$ws.Range("N258").Value = 6000
$ws.Range("O258").Value = 6000
$ws.Range("P258").Value = 6000
$ws.Range("Q258").Value = "$/bandeja 7 kilos"
$ws.Range("R258").Value = "Provincia de Melipilla"
$ws.Range("S258").Value = 857
$ws.Range("T258").Value = 7

# --- New row 259 ---
$ws.Range("A259").Value = 3
$ws.Range("B259").Value = "Femacal de La Calera"
$ws.Range("C259").Value = "Coquimbo"
$ws.Range("D259").Value = 44875
$ws.Range("E259").Value = 5
$ws.Range("F259").Value = "Fruta"
$ws.Range("G259").Value = 100101
$ws.Range("H259").Value = "Berries"
$ws.Range("I259").Value = 100112025
$ws.Range("J259").Value = "Frutilla"
$ws.Range("K259").Value = "Sin especificar"
$ws.Range("L259").Value = "Segunda"
$ws.Range("M259").Value = 56
$ws.Range("N259").Value = 4000
$ws.Range("O259").Value = 4000
$ws.Range("P259").Value = 4000
$ws.Range("Q259").Value = "$/bandeja 7 kilos"
$ws.Range("R259").Value = "Provincia de Melipilla"
$ws.Range("S259").Value = 571
$ws.Range("T259").Value = 7

Write-Output "OK"
